$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E4").Value = 1184000.000000039
$ws.Range("F6").Value = 1177703.294781314
$ws.Range("G6").Value = 1155650.958798337
$ws.Range("G7").Value = 22047.71408006548
